$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "90.761.67"
$ws.Range("E2").Value = "  +1.31%  "

# Row 3
Set-TextValue "D3" "3.161.32"
$ws.Range("E3").Value = "  +2.89%  "

# Row 4
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
Set-TextValue "D5" "217.45"
$ws.Range("E5").Value = "  +2.07%  "

# Row 6
Set-TextValue "D6" "626.71"
$ws.Range("E6").Value = "  +2.10%  "

# Row 7
Set-TextValue "D7" "1.15"
$ws.Range("E7").Value = "  +31.37%  "

# Row 8
Set-TextValue "D8" "0.368"
$ws.Range("E8").Value = "  +0.71%  "

# Row 9
Set-TextValue "D9" "1.00"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
Set-TextValue "D10" "3.158.70"
$ws.Range("E10").Value = "  +2.85%  "

# Row 11
Set-TextValue "D11" "0.762"
$ws.Range("E11").Value = "  +15.24%  "

# Row 12
Set-TextValue "D12" "0.201"
$ws.Range("E12").Value = "  +6.78%  "

# Row 13
Set-TextValue "D13" "0.0000246"
$ws.Range("E13").Value = "  +2.96%  "

# Row 14
Set-TextValue "D14" "5.70"
$ws.Range("E14").Value = "  +6.52%  "

# Row 15
Set-TextValue "D15" "35.19"
$ws.Range("E15").Value = "  +8.91%  "

# Row 16
Set-TextValue "D16" "90.715.36"
$ws.Range("E16").Value = "  +1.33%  "

# Row 17
Set-TextValue "D17" "3.748.57"
$ws.Range("E17").Value = "  +2.40%  "

# Row 18
Set-TextValue "D18" "3.113.23"
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
Set-TextValue "D19" "3.71"
$ws.Range("E19").Value = "  +9.30%  "

# Row 20
Set-TextValue "D20" "0.0000215"
$ws.Range("E20").Value = "  +2.31%  "

# Row 21
Set-TextValue "D21" "14.29"
$ws.Range("E21").Value = "  +5.78%  "

# Row 22
Set-TextValue "D22" "448.71"
$ws.Range("E22").Value = "  +4.19%  "

# Row 23
Set-TextValue "D23" "8.98"
$ws.Range("E23").Value = "  +9.77%  "

# Row 24
Set-TextValue "D24" "5.22"
$ws.Range("E24").Value = "  +4.44%  "

# Row 25
Set-TextValue "D25" "5.98"
$ws.Range("E25").Value = "  +9.35%  "

# Row 26
Set-TextValue "D26" "89.45"
$ws.Range("E26").Value = "  +5.14%  "

# Row 27
Set-TextValue "D27" "12.19"
$ws.Range("E27").Value = "  +1.43%  "

# Row 28
Set-TextValue "D28" "3.336.85"
$ws.Range("E28").Value = "  +2.32%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "9.29"
$ws.Range("E30").Value = "  +14.04%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D31" "0.162"
$ws.Range("E31").Value = "  -1.52%  "

# Row 32
Set-TextValue "D32" "1.02"
$ws.Range("E32").Value = "  -7.35%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "527.63"
$ws.Range("E33").Value = "  +4.15%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "25.85"
$ws.Range("E34").Value = "  +13.63%  "

# Row 35
Set-TextValue "D35" "3.68"
$ws.Range("E35").Value = "  +1.88%  "

# Row 36
$ws.Range("E36").Value = "  +9.50%  "

# Row 37
Set-TextValue "D37" "6.98"
$ws.Range("E37").Value = "  +5.08%  "

# Row 38
$ws.Range("E38").Value = "  +6.64%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.176"
$ws.Range("E39").Value = "  +27.07%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D40" "1.30"
$ws.Range("E40").Value = "  +4.45%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0871"
$ws.Range("E41").Value = "  +25.92%  "

# Row 42
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D42" "22.21"
$ws.Range("E42").Value = "  -0.33%  "

# Row 43
$ws.Range("E43").Value = "  -0.24%  "

# Row 44
Set-TextValue "D44" "0.418"
$ws.Range("E44").Value = "  +13.09%  "

# Row 45
$ws.Range("E45").Value = "  +6.22%  "

# Row 46
$ws.Range("E46").Value = "  -0.02%  "

# Row 47
Set-TextValue "D47" "147.42"
$ws.Range("E47").Value = "  +0.07%  "

# Row 48
Set-TextValue "D48" "1.34"
$ws.Range("E48").Value = "  +10.43%  "

# Row 49
Set-TextValue "D49" "44.26"
$ws.Range("E49").Value = "  +1.13%  "

# Row 50
Set-TextValue "D50" "4.41"
$ws.Range("E50").Value = "  +8.57%  "

# Row 51
Set-TextValue "D51" "0.652"
$ws.Range("E51").Value = "  +11.08%  "
